$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Books")

for ($r = 5; $r -le 256; $r++) {
    $ws.Cells.Item($r, 20).Value = 0
}

$ws.Activate()
$ws.Range("T259").Select()
Write-Output "done"
